$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Paragraph 1 ends with "This is the first file" and, right before its
# paragraph mark, carries a collapsed "_GoBack" bookmark. We need to:
#   1. Split off a new, second paragraph after that text.
#   2. Fill the new paragraph with "This is 2" + "nd" (superscript) + " commit"
#   3. Leave the _GoBack bookmark collapsed at the end of the new
#      paragraph (i.e. it should move along with the edit).
# ------------------------------------------------------------------

$firstPara = $d.Paragraphs(1).Range
$splitPoint = $firstPara.End - 1   # just before paragraph 1's own paragraph mark

$ip = $d.Range($splitPoint, $splitPoint)
$ip.InsertParagraphAfter() | Out-Null

# The old bookmark remains attached to paragraph 1; drop it, it will be
# rebuilt at the correct location once the new paragraph has its text.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Type the new paragraph's full text first as plain text so later
# position-based edits (bookmark placement) operate on a single,
# un-split run.
$newParaRange = $d.Paragraphs(2).Range
$newParaRange.Collapse(1)   # wdCollapseStart
$newParaRange.InsertAfter("This is 2nd commit")
$newParaRange.Collapse(0)   # wdCollapseEnd

# ------------------------------------------------------------------
# Recreate the _GoBack bookmark collapsed at the very end of the new
# paragraph. Adding a bookmark directly at a paragraph-end boundary is
# unreliable, so add it one character earlier (a safe position) and
# then nudge it forward by retyping the final character -- retyping
# text exactly at a collapsed bookmark's position naturally extends
# the bookmark across the retyped text.
# ------------------------------------------------------------------

$endPos = $newParaRange.End
$anchorPos = $endPos - 1

$bookmarkSeed = $d.Range($anchorPos, $anchorPos)
$d.Bookmarks.Add("_GoBack", $bookmarkSeed)

$lastCharRange = $d.Range($anchorPos, $endPos)
$lastChar = $lastCharRange.Text
$lastCharRange.Delete()

$retypeIp = $d.Range($anchorPos, $anchorPos)
$retypeIp.InsertAfter($lastChar)

# ------------------------------------------------------------------
# Now that the text and bookmark are in place, apply superscript
# formatting to "nd" within the new paragraph.
# ------------------------------------------------------------------

$p2 = $d.Paragraphs(2).Range
$ndStart = $p2.Start + ("This is 2").Length
$ndEnd = $ndStart + ("nd").Length
$ndRange = $d.Range($ndStart, $ndEnd)
$ndRange.Font.Superscript = $true
